$wb = $excel.ActiveWorkbook

# The data lives on the "Global" worksheet (first sheet in the workbook).
$wsGlobal = $wb.Worksheets.Item("Global")
# "GLOverview" is the sheet that is actually active/selected in the workbook.
$wsOverview = $wb.Worksheets.Item("GLOverview")

# Update the browser name used for the test from CHROME to FIREFOX.
$wsGlobal.Range("A2").Value = "FIREFOX"

# Move the selection on the "Global" sheet to A2 (it was B2 before the edit),
# then re-activate the "GLOverview" sheet so the workbook's active tab is
# left unchanged.
$wsGlobal.Range("A2").Select()
$wsOverview.Activate()
